$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data cells in this sheet are stored as literal text, even the
# numeric-looking ones (row index, barcodes with leading zeros, etc.) - that
# is why the sheet already carries a "number stored as text" ignored-error
# hint. Force Text format on every cell we are about to (re)write so values
# that look numeric (e.g. "3", "0000000000001") are kept as text instead of
# being coerced to numbers. A2:A4 are left untouched since their values
# ("0","1","2") are not changing.
$ws.Range("C2:D2").NumberFormat = "@"
$ws.Range("C3:D3").NumberFormat = "@"
$ws.Range("C4:D4").NumberFormat = "@"
$ws.Range("A5:D11").NumberFormat = "@"

# Existing rows 2-4: product/type columns change, barcode column changes.
$ws.Range("B2").Value = "plastic"
$ws.Range("C2").Value = "0000000000001"
$ws.Range("D2").Value = "plastic"

$ws.Range("B3").Value = "paper"
$ws.Range("C3").Value = "0000000000002"
$ws.Range("D3").Value = "paper"

$ws.Range("B4").Value = "glass"
$ws.Range("C4").Value = "0000000000003"
$ws.Range("D4").Value = "glass"

# Existing row 5: full row replaced with new product.
$ws.Range("A5").Value = "3"
$ws.Range("B5").Value = "Calve knoflook saus"
$ws.Range("C5").Value = "8720182255563"
$ws.Range("D5").Value = "plastic"

# New rows 6-11.
$ws.Range("A6").Value = "4"
$ws.Range("B6").Value = "Verstegen paprikapoeder"
$ws.Range("C6").Value = "8712200856104"
$ws.Range("D6").Value = "glass"

$ws.Range("A7").Value = "5"
$ws.Range("B7").Value = "AH tomatenpuree"
$ws.Range("C7").Value = "8059602910011"
$ws.Range("D7").Value = "plastic"

$ws.Range("A8").Value = "6"
$ws.Range("B8").Value = "Calve pindakaas stukjes"
$ws.Range("C8").Value = "8711200430925"
$ws.Range("D8").Value = "glass"

$ws.Range("A9").Value = "7"
$ws.Range("B9").Value = "Cup a soup tomaat"
$ws.Range("C9").Value = "5711327460348"
$ws.Range("D9").Value = "glass"

$ws.Range("A10").Value = "8"
$ws.Range("B10").Value = "Pickwick winterglow"
$ws.Range("C10").Value = "8711000008881"
$ws.Range("D10").Value = "paper"

$ws.Range("A11").Value = "9"
$ws.Range("B11").Value = "Coca Cola blik"
$ws.Range("C11").Value = "5449000008046"
$ws.Range("D11").Value = "paper"
